$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1164.6666
$ws.Range("J12").Value = 5263
$ws.Range("L12").Value = 5263
$ws.Range("N12").Value = -5603
$ws.Range("H15").Value = 3713205
$ws.Range("I15").Value = 3713205
$ws.Range("K15").Value = 11139615
$ws.Range("M15").Value = -11139446
$ws.Range("H33").Value = 274.57895
$ws.Range("I33").Value = 250.71428
$ws.Range("K33").Value = 250.71428
$ws.Range("M33").Value = -21.71428
$ws.Range("H51").Value = 22699.934
$ws.Range("I51").Value = 10812.375
$ws.Range("K51").Value = 10812.375
$ws.Range("M51").Value = -10328.375
$ws.Range("H112").Value = 3342.3542
$ws.Range("J112").Value = 3342.3542
$ws.Range("L112").Value = 10027.0626
$ws.Range("N112").Value = -12243.0626
$ws.Range("H137").Value = 6584394.5
$ws.Range("I137").Value = 12503833
$ws.Range("K137").Value = 37511499
$ws.Range("M137").Value = -37508949
$ws.Range("H138").Value = 2628.205
$ws.Range("I138").Value = 1610.3125
$ws.Range("J138").Value = 3336.3044
$ws.Range("K138").Value = 4830.9375
$ws.Range("L138").Value = 10008.9132
$ws.Range("M138").Value = 309.0625
$ws.Range("N138").Value = -20288.9132
$ws.Range("H141").Value = 1133.8
$ws.Range("I141").Value = 1000
$ws.Range("K141").Value = 3000
$ws.Range("M141").Value = 2180

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 34861
$ws.Range("I32").Value = 34861
$ws.Range("K32").Value = 34861
$ws.Range("M32").Value = -34574
$ws.Range("H34").Value = 289666.66
$ws.Range("J34").Value = 500000
$ws.Range("L34").Value = 500000
$ws.Range("N34").Value = -500542
$ws.Range("H102").Value = 3260.5557
$ws.Range("I102").Value = 3260.5557
$ws.Range("K102").Value = 3260.5557
$ws.Range("M102").Value = -1638.5557
$ws.Range("H132").Value = 5961.625
$ws.Range("I132").Value = 4787.7
$ws.Range("J132").Value = 7918.1665
$ws.Range("K132").Value = 14363.1
$ws.Range("L132").Value = 23754.4995
$ws.Range("M132").Value = -11833.1
$ws.Range("N132").Value = -28814.4995
$ws.Range("H141").Value = 109088.664
$ws.Range("J141").Value = 109088.664
$ws.Range("L141").Value = 109088.664
$ws.Range("N141").Value = -119448.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 168854.67
$ws.Range("I86").Value = 2174.9
$ws.Range("J86").Value = 1002253.5
$ws.Range("K86").Value = 2174.9
$ws.Range("L86").Value = 1002253.5
$ws.Range("M86").Value = -1051.9
$ws.Range("N86").Value = -1004499.5
$ws.Range("H89").Value = 168854.67
$ws.Range("I89").Value = 2174.9
$ws.Range("J89").Value = 1002253.5
$ws.Range("K89").Value = 10874.5
$ws.Range("L89").Value = 5011267.5
$ws.Range("M89").Value = -5258.5
$ws.Range("N89").Value = -5022499.5
$ws.Range("H99").Value = 2970.6086
$ws.Range("I99").Value = 2322.1333
$ws.Range("K99").Value = 2322.1333
$ws.Range("M99").Value = -824.1333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1871.0526
$ws.Range("I22").Value = 763.5
$ws.Range("J22").Value = 3101.6667
$ws.Range("K22").Value = 763.5
$ws.Range("L22").Value = 3101.6667
$ws.Range("M22").Value = -413.5
$ws.Range("N22").Value = -3801.6667
$ws.Range("H31").Value = 37040250
$ws.Range("I31").Value = 52632892
$ws.Range("K31").Value = 52632892
$ws.Range("M31").Value = -52632597
$ws.Range("H33").Value = 1633
$ws.Range("I33").Value = 1633
$ws.Range("K33").Value = 1633
$ws.Range("M33").Value = -1254
$ws.Range("H34").Value = 37040250
$ws.Range("I34").Value = 52632892
$ws.Range("K34").Value = 52632892
$ws.Range("M34").Value = -52632690
$ws.Range("H58").Value = 5164
$ws.Range("I58").Value = 3209.1428
$ws.Range("K58").Value = 3209.1428
$ws.Range("M58").Value = -3006.1428
$ws.Range("H74").Value = 38297.75
$ws.Range("J74").Value = 38297.75
$ws.Range("L74").Value = 38297.75
$ws.Range("N74").Value = -40045.75
$ws.Range("H77").Value = 38297.75
$ws.Range("J77").Value = 38297.75
$ws.Range("L77").Value = 114893.25
$ws.Range("N77").Value = -123629.25
$ws.Range("H132").Value = 206000.17
$ws.Range("I132").Value = 4625.25
$ws.Range("K132").Value = 13875.75
$ws.Range("M132").Value = -11345.75
$ws.Range("H134").Value = 3644.913
$ws.Range("I134").Value = 3339.322
$ws.Range("J134").Value = 5447.9
$ws.Range("K134").Value = 10017.966
$ws.Range("L134").Value = 16343.7
$ws.Range("M134").Value = -7482.966
$ws.Range("N134").Value = -21413.7
$ws.Range("H136").Value = 5164
$ws.Range("I136").Value = 3209.1428
$ws.Range("K136").Value = 9627.428400000001
$ws.Range("M136").Value = -7077.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4153691.8
$ws.Range("I4").Value = 3273987.8
$ws.Range("K4").Value = 9821963.399999999
$ws.Range("M4").Value = -9821851.399999999
$ws.Range("H39").Value = 2254
$ws.Range("J39").Value = 3600
$ws.Range("L39").Value = 10800
$ws.Range("N39").Value = -11388
$ws.Range("H131").Value = 23819872
$ws.Range("I131").Value = 83334230
$ws.Range("J131").Value = 14127
$ws.Range("K131").Value = 250002690
$ws.Range("L131").Value = 42381
$ws.Range("M131").Value = -249997650
$ws.Range("N131").Value = -52461

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4991.6665
$ws.Range("I80").Value = 4479
$ws.Range("K80").Value = 4479
$ws.Range("M80").Value = -3481
$ws.Range("H83").Value = 4991.6665
$ws.Range("I83").Value = 4479
$ws.Range("K83").Value = 22395
$ws.Range("M83").Value = -17403
$ws.Range("H122").Value = 6943.1924
$ws.Range("I122").Value = 8370.9
$ws.Range("J122").Value = 2184.1667
$ws.Range("K122").Value = 25112.7
$ws.Range("L122").Value = 6552.500100000001
$ws.Range("M122").Value = -22662.7
$ws.Range("N122").Value = -11452.5001
$ws.Range("H132").Value = 6019.625
$ws.Range("I132").Value = 3119.4546
$ws.Range("J132").Value = 12400
$ws.Range("K132").Value = 9358.363799999999
$ws.Range("L132").Value = 37200
$ws.Range("M132").Value = -6828.363799999999
$ws.Range("N132").Value = -42260

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3206.2432
$ws.Range("J22").Value = 4456.3887
$ws.Range("L22").Value = 4456.3887
$ws.Range("N22").Value = -5046.3887
$ws.Range("H27").Value = 3206.2432
$ws.Range("J27").Value = 4456.3887
$ws.Range("L27").Value = 4456.3887
$ws.Range("N27").Value = -4670.3887
$ws.Range("H32").Value = 1945.125
$ws.Range("I32").Value = 1945.125
$ws.Range("K32").Value = 1945.125
$ws.Range("M32").Value = -1628.125
$ws.Range("H46").Value = 6381.6772
$ws.Range("I46").Value = 2104.7144
$ws.Range("J46").Value = 7629.125
$ws.Range("K46").Value = 2104.7144
$ws.Range("L46").Value = 7629.125
$ws.Range("M46").Value = -1916.7144
$ws.Range("N46").Value = -8005.125
$ws.Range("H55").Value = 29972.059
$ws.Range("I55").Value = 50668.7
$ws.Range("K55").Value = 50668.7
$ws.Range("M55").Value = -50495.7
$ws.Range("H93").Value = 2138.6667
$ws.Range("I93").Value = 1257.3334
$ws.Range("K93").Value = 1257.3334
$ws.Range("M93").Value = -9.333399999999983
$ws.Range("H100").Value = 16669026
$ws.Range("I100").Value = 41668264
$ws.Range("J100").Value = 2865.3333
$ws.Range("K100").Value = 41668264
$ws.Range("L100").Value = 2865.3333
$ws.Range("M100").Value = -41667723
$ws.Range("N100").Value = -3947.3333
$ws.Range("H132").Value = 12472.5
$ws.Range("I132").Value = 7495
$ws.Range("J132").Value = 17450
$ws.Range("K132").Value = 22485
$ws.Range("L132").Value = 52350
$ws.Range("M132").Value = -19955
$ws.Range("N132").Value = -57410
$ws.Range("H136").Value = 5629.9546
$ws.Range("I136").Value = 3500.2856
$ws.Range("K136").Value = 10500.8568
$ws.Range("M136").Value = -7950.856800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 9153.143
$ws.Range("I58").Value = 9012
$ws.Range("K58").Value = 9012
$ws.Range("M58").Value = -8704
$ws.Range("H81").Value = 14785.517
$ws.Range("I81").Value = 3698.75
$ws.Range("J81").Value = 16428
$ws.Range("K81").Value = 7397.5
$ws.Range("L81").Value = 32856
$ws.Range("M81").Value = -6336.5
$ws.Range("N81").Value = -34978
$ws.Range("H84").Value = 14785.517
$ws.Range("I84").Value = 3698.75
$ws.Range("J84").Value = 16428
$ws.Range("K84").Value = 36987.5
$ws.Range("L84").Value = 164280
$ws.Range("M84").Value = -31683.5
$ws.Range("N84").Value = -174888
$ws.Range("H100").Value = 5599
$ws.Range("I100").Value = 5983.9
$ws.Range("J100").Value = 1750
$ws.Range("K100").Value = 11967.8
$ws.Range("L100").Value = 3500
$ws.Range("M100").Value = -11426.8
$ws.Range("N100").Value = -4582
$ws.Range("H107").Value = 4438.8335
$ws.Range("I107").Value = 4499.857
$ws.Range("J107").Value = 4353.4
$ws.Range("K107").Value = 13499.571
$ws.Range("L107").Value = 13060.2
$ws.Range("M107").Value = -11579.571
$ws.Range("N107").Value = -16900.2
$ws.Range("H122").Value = 3986.111
$ws.Range("I122").Value = 3986.111
$ws.Range("K122").Value = 11958.333
$ws.Range("M122").Value = -9508.332999999999
$ws.Range("H132").Value = 7193.3335
$ws.Range("I132").Value = 5333.3335
$ws.Range("K132").Value = 16000.0005
$ws.Range("M132").Value = -13470.0005
$ws.Range("H136").Value = 5098.222
$ws.Range("I136").Value = 1477.6
$ws.Range("K136").Value = 4432.799999999999
$ws.Range("M136").Value = -1882.799999999999
